$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C ("Förändrad") for rows 2 through 11 from 45224 to 45233.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
